# Update Efna5-Epha4.xlsx (LR-pairs_lrc2p) with refreshed TPM-derived NATMI
# values: a third sending/target cluster ("ECs") is now present, so the
# original 2x3 = 6 result rows become 3x3 = 9 result rows, and every
# specificity-style column is recomputed against the enlarged cluster set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (columns A:T), in final row order:
# Sending cluster cycles ECs -> FAPs -> MuSCs, Target cluster cycles
# ECs -> FAPs -> MuSCs within each sending block.
$rows = @(
  @("ECs",   "Efna5", "Epha4", "ECs",   1, 0.3333333333333333, 0.05800433333333333, 0.174013,  0.02087975181349295, 0.02087975181349295, 3, 1, 8.081040666666667, 24.243122, 0.4661250698616886, 0.4661250698616886, 0.4687353765095555,  4.218618388586,     0.009732575772759123, 0.009732575772759123),
  @("ECs",   "Efna5", "Epha4", "FAPs",  1, 0.3333333333333333, 0.05800433333333333, 0.174013,  0.02087975181349295, 0.02087975181349295, 3, 1, 7.000300666666668, 21.000902, 0.4037865631294714, 0.4037865631294715, 0.4060477733028889,  3.654429959726001,  0.008430963223766668, 0.008430963223766668),
  @("ECs",   "Efna5", "Epha4", "MuSCs", 1, 0.3333333333333333, 0.05800433333333333, 0.174013,  0.02087975181349295, 0.02087975181349295, 3, 1, 2.255294666666666, 6.765884,  0.1300883670088399, 0.1300883670088399, 0.1308168636102222,  1.177351772492,     0.002716212816967162, 0.002716212816967162),
  @("FAPs",  "Efna5", "Epha4", "ECs",   3, 1,                  1.666083666666667,  4.998251,  0.5997381826733804,  0.5997381826733805,  3, 1, 8.081040666666667, 24.243122, 0.4661250698616886, 0.4661250698616886, 13.46368986440244,   121.173208779622,   0.2795530022973516,  0.2795530022973516),
  @("FAPs",  "Efna5", "Epha4", "FAPs",  3, 1,                  1.666083666666667,  4.998251,  0.5997381826733804,  0.5997381826733805,  3, 1, 7.000300666666668, 21.000902, 0.4037865631294714, 0.4037865631294715, 11.66308660248911,   104.967779422402,   0.2421662195591994,  0.2421662195591995),
  @("FAPs",  "Efna5", "Epha4", "MuSCs", 3, 1,                  1.666083666666667,  4.998251,  0.5997381826733804,  0.5997381826733805,  3, 1, 2.255294666666666, 6.765884,  0.1300883670088399, 0.1300883670088399, 3.757509607653777,   33.817586468884,    0.07801896081682938, 0.07801896081682939),
  @("MuSCs", "Efna5", "Epha4", "ECs",   3, 1,                  1.053930333333333,  3.161791,  0.3793820655131266,  0.3793820655131266,  3, 1, 8.081040666666667, 24.243122, 0.4661250698616886, 0.4661250698616886, 8.516853883500223,   76.65168495150201,  0.1768394917915778,  0.1768394917915778),
  @("MuSCs", "Efna5", "Epha4", "FAPs",  3, 1,                  1.053930333333333,  3.161791,  0.3793820655131266,  0.3793820655131266,  3, 1, 7.000300666666668, 21.000902, 0.4037865631294714, 0.4037865631294715, 7.377829215053557,   66.40046293548201,  0.1531893803465053,  0.1531893803465054),
  @("MuSCs", "Efna5", "Epha4", "MuSCs", 3, 1,                  1.053930333333333,  3.161791,  0.3793820655131266,  0.3793820655131266,  3, 1, 2.255294666666666, 6.765884,  0.1300883670088399, 0.1300883670088399, 2.376923459804889,   21.392311138244,    0.04935319337504336, 0.04935319337504336)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $startRow + $i
  $rowData = $rows[$i]
  for ($c = 1; $c -le $rowData.Count; $c++) {
    $ws.Cells.Item($r, $c).Value2 = $rowData[$c - 1]
  }
}

Write-Host "Updated rows 2-10 with refreshed TPM NATMI values (added ECs cluster)."
